# Apply updated "想去人数" (want-to-go count) values to the 展览 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7288
$ws1.Range("F4").Value = 467
$ws1.Range("F7").Value = 196
$ws1.Range("F12").Value = 225
$ws1.Range("F19").Value = 3815
$ws1.Range("F26").Value = 2498
$ws1.Range("F28").Value = 325
$ws1.Range("F33").Value = 27
$ws1.Range("F38").Value = 35
$ws1.Range("F39").Value = 1491

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7288
$ws4.Range("F4").Value = 467
$ws4.Range("F8").Value = 196
$ws4.Range("F13").Value = 225
$ws4.Range("F20").Value = 3815
$ws4.Range("F27").Value = 2498
$ws4.Range("F29").Value = 325
$ws4.Range("F34").Value = 27
$ws4.Range("F39").Value = 35
$ws4.Range("F40").Value = 1491
